# Refresh market-price snapshot data (currentAveragePrice* / Leve profit
# columns H:N) across all 8 profession sheets to match the latest pull from
# the scheduled market-data runner. Pure data refresh - no formulas, no
# structural changes. A few rows also gain/lose a trailing HQ-profit cell
# (column N or M) depending on whether an HQ recipe/listing now exists.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2458.6667
$ws.Range("I17").Value = 900
$ws.Range("J17").Value = 2600.3635
$ws.Range("K17").Value = 2700
$ws.Range("L17").Value = 7801.0905
$ws.Range("M17").Value = -2532
$ws.Range("N17").Value = -8137.0905
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H29").Value = 47
$ws.Range("I29").Value = 47
$ws.Range("K29").Value = 141
$ws.Range("M29").Value = 140
$ws.Range("H31").Value = 209.28572
$ws.Range("I31").Value = 44.166668
$ws.Range("K31").Value = 132.500004
$ws.Range("M31").Value = 97.49999600000001
$ws.Range("H33").Value = 104.333336
$ws.Range("I33").Value = 85.2
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 85.2
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = 143.8
$ws.Range("N33").Value = -658
$ws.Range("H69").Value = 7994.3335
$ws.Range("I69").Value = 8218.666999999999
$ws.Range("K69").Value = 24656.001
$ws.Range("M69").Value = -23782.001
$ws.Range("H72").Value = 7994.3335
$ws.Range("I72").Value = 8218.666999999999
$ws.Range("K72").Value = 73968.003
$ws.Range("M72").Value = -69600.003
$ws.Range("H76").Value = 6833.1665
$ws.Range("I76").Value = 5750
$ws.Range("K76").Value = 5750
$ws.Range("M76").Value = -5435
$ws.Range("H79").Value = 6833.1665
$ws.Range("I79").Value = 5750
$ws.Range("K79").Value = 5750
$ws.Range("M79").Value = -4658
$ws.Range("H94").Value = 3129.2856
$ws.Range("I94").Value = 3129.2856
$ws.Range("K94").Value = 3129.2856
$ws.Range("M94").Value = -2678.2856
$ws.Range("H116").Value = 4499
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2636031.8
$ws.Range("I32").Value = 2770.7058
$ws.Range("K32").Value = 2770.7058
$ws.Range("M32").Value = -2483.7058
$ws.Range("H45").Value = 3915
$ws.Range("I45").Value = 2163.3333
$ws.Range("K45").Value = 2163.3333
$ws.Range("M45").Value = -1786.3333
$ws.Range("H88").Value = 1243.8889
$ws.Range("I88").Value = 595
$ws.Range("K88").Value = 595
$ws.Range("M88").Value = -189
$ws.Range("H91").Value = 1243.8889
$ws.Range("I91").Value = 595
$ws.Range("K91").Value = 595
$ws.Range("M91").Value = 809
$ws.Range("H97").Value = 1221.4286
$ws.Range("I97").Value = 887.9
$ws.Range("K97").Value = 887.9
$ws.Range("M97").Value = -391.9
$ws.Range("H102").Value = 4161.727
$ws.Range("I102").Value = 825.5714
$ws.Range("K102").Value = 825.5714
$ws.Range("M102").Value = 796.4286
$ws.Range("H110").Value = 2606.3635
$ws.Range("I110").Value = 2327.5
$ws.Range("J110").Value = 2941
$ws.Range("K110").Value = 2327.5
$ws.Range("L110").Value = 2941
$ws.Range("M110").Value = -282.5
$ws.Range("N110").Value = -7031
$ws.Range("H122").Value = 1694.9166
$ws.Range("I122").Value = 990.5
$ws.Range("J122").Value = 3103.75
$ws.Range("K122").Value = 2971.5
$ws.Range("L122").Value = 9311.25
$ws.Range("M122").Value = -521.5
$ws.Range("N122").Value = -14211.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 866.6667
$ws.Range("I94").Value = 1000
$ws.Range("K94").Value = 1000
$ws.Range("M94").Value = -549
$ws.Range("H99").Value = 1652.8
$ws.Range("I99").Value = 1666
$ws.Range("K99").Value = 1666
$ws.Range("M99").Value = -168
$ws.Range("H105").Value = 1741
$ws.Range("I105").Value = 1764.5
$ws.Range("J105").Value = 1600
$ws.Range("K105").Value = 1764.5
$ws.Range("L105").Value = 1600
$ws.Range("M105").Value = -17.5
$ws.Range("N105").Value = -5094
$ws.Range("H107").Value = 7044
$ws.Range("I107").Value = 1831.3334
$ws.Range("J107").Value = 8998.75
$ws.Range("K107").Value = 1831.3334
$ws.Range("L107").Value = 8998.75
$ws.Range("M107").Value = 88.66660000000002
$ws.Range("N107").Value = -12838.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1342.2
$ws.Range("I16").Value = 1252.875
$ws.Range("K16").Value = 1252.875
$ws.Range("M16").Value = -965.875
$ws.Range("H62").Value = 5958.3335
$ws.Range("J62").Value = 5958.3335
$ws.Range("L62").Value = 5958.3335
$ws.Range("N62").Value = -7206.3335
$ws.Range("H65").Value = 5958.3335
$ws.Range("J65").Value = 5958.3335
$ws.Range("L65").Value = 29791.6675
$ws.Range("N65").Value = -36031.6675
$ws.Range("H107").Value = 1567.6
$ws.Range("I107").Value = 1352.6666
$ws.Range("K107").Value = 1352.6666
$ws.Range("M107").Value = 567.3334
$ws.Range("H113").Value = 1342.2
$ws.Range("I113").Value = 1252.875
$ws.Range("K113").Value = 1252.875
$ws.Range("M113").Value = 917.125
$ws.Range("H132").Value = 4259.8
$ws.Range("I132").Value = 4199.75
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 12599.25
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -10069.25
$ws.Range("N132").Value = -18560

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 4555.3335
$ws.Range("J19").Value = 4555.3335
$ws.Range("L19").Value = 13666.0005
$ws.Range("N19").Value = -14014.0005
$ws.Range("H107").Value = 744
$ws.Range("I107").Value = 488
$ws.Range("K107").Value = 1464
$ws.Range("M107").Value = 456
$ws.Range("H121").Value = 933.625
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 995.5714
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 2986.7142
$ws.Range("M121").Value = -190
$ws.Range("N121").Value = -5606.7142
$ws.Range("H141").Value = 1175
$ws.Range("I141").Value = 1175
$ws.Range("K141").Value = 3525
$ws.Range("M141").Value = 1655

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5313.2
$ws.Range("I80").Value = 6200
$ws.Range("J80").Value = 4722
$ws.Range("K80").Value = 6200
$ws.Range("L80").Value = 4722
$ws.Range("M80").Value = -5202
$ws.Range("N80").Value = -6718
$ws.Range("H83").Value = 5313.2
$ws.Range("I83").Value = 6200
$ws.Range("J83").Value = 4722
$ws.Range("K83").Value = 31000
$ws.Range("L83").Value = 23610
$ws.Range("M83").Value = -26008
$ws.Range("N83").Value = -33594
$ws.Range("H97").Value = 1299.6666
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1299.6666
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 1299.6666
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -2291.6666
$ws.Range("H99").Value = 3607.3333
$ws.Range("I99").Value = 2058.25
$ws.Range("K99").Value = 2058.25
$ws.Range("M99").Value = 187.75
$ws.Range("H122").Value = 3161.625
$ws.Range("I122").Value = 2553.6924
$ws.Range("K122").Value = 7661.0772
$ws.Range("M122").Value = -5211.0772
$ws.Range("H126").Value = 3402.7058
$ws.Range("I126").Value = 2507.8
$ws.Range("J126").Value = 4681.143
$ws.Range("K126").Value = 7523.400000000001
$ws.Range("L126").Value = 14043.429
$ws.Range("M126").Value = -5053.400000000001
$ws.Range("N126").Value = -18983.429

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 100
$ws.Range("J2").Value = 100
$ws.Range("L2").Value = 100
$ws.Range("N2").Value = -324
$ws.Range("H61").Value = 3930.4375
$ws.Range("I61").Value = 2088.9
$ws.Range("K61").Value = 2088.9
$ws.Range("M61").Value = -1886.9
$ws.Range("H93").Value = 1688.4
$ws.Range("I93").Value = 1610.5
$ws.Range("K93").Value = 1610.5
$ws.Range("M93").Value = -362.5
$ws.Range("H113").Value = 3930.4375
$ws.Range("I113").Value = 2088.9
$ws.Range("K113").Value = 2088.9
$ws.Range("M113").Value = 81.09999999999991
$ws.Range("H122").Value = 2725.75
$ws.Range("I122").Value = 2725.75
$ws.Range("K122").Value = 8177.25
$ws.Range("M122").Value = -5727.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1370
$ws.Range("I100").Value = 1370
$ws.Range("K100").Value = 2740
$ws.Range("M100").Value = -2199
